$wb = $excel.ActiveWorkbook

# Sheet "2" (first sheet): the row 10 scratch entries (B10 "PevalAll" /
# C10 the recursive eval_all pyeval formula) are removed; the
# eval_all formula is relocated onto a brand-new "e2" sheet (see below),
# and "PevalAll" is dropped entirely.
$ws1 = $wb.Worksheets.Item("2")
$ws1.Range("B10:C10").ClearContents()

# "eval sheet": add a new D4 cell with a short cross-sheet pointer to the
# new "e2" sheet's B2 ("foo") cell.
$wsEval = $wb.Worksheets.Item("eval sheet")
$wsEval.Range("D4").Value = "#B2"

# Add a brand-new worksheet named "e2" as the last tab, holding the
# recursive pyeval formula that used to live at "2"!C10.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsE2 = $wb.Worksheets.Add($null, $lastSheet)
$wsE2.Name = "e2"
$wsE2.Range("A1").Value = '#eval sheet!::{"func": "pipe", "kwds":{"lax": false}, "args":[["df", {"index_col": null}], ["pyeval", {"include": "EVAL_COL", "eval_all": true}], "recurse"]}'

# Restore sheet "2" as the active sheet, with the selection moved from the
# deleted C10 onto B10.
$ws1.Activate()
$ws1.Range("B10").Select()
